# Update the "Estado de Cuenta" worker rows (16-21):
#  - Periodo Mora (column E) is reordered from descending (1903..1810)
#    to ascending (1810..1903) chronological order, keeping each
#    period's own "Valor Mora" (column F) value attached to it.
#  - Salario Basico (column G) is updated to the new value 737717
#    for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periodos = @("1810", "1811", "1812", "1901", "1902", "1903")
$valoresMora = @(29509, 29509, 29509, 29509, 29509, 24591)
$salarioBasico = 737717

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valoresMora[$i]
    $ws.Cells.Item($row, 7).Value = $salarioBasico
}
